# Updated: po 07. 03. 2022
# Refresh the OpenData Slovakia Covid daily-stats sheet:
#  - revise AgTests (F) / AgPosit (G) cumulative figures for 2022-02-25 .. 2022-03-02
#  - append four new daily rows for 2022-03-03 .. 2022-03-06

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Revise existing AgTests / AgPosit figures (rows 697-728) ---
$ws.Range("F697").Value = 29067
$ws.Range("G697").Value = 3055
$ws.Range("F698").Value = 70829
$ws.Range("F699").Value = 43577
$ws.Range("F701").Value = 41857
$ws.Range("F702").Value = 36324
$ws.Range("G702").Value = 3928
$ws.Range("F703").Value = 17081
$ws.Range("G703").Value = 2609
$ws.Range("F704").Value = 25057
$ws.Range("G704").Value = 3706
$ws.Range("F705").Value = 55996
$ws.Range("G705").Value = 6307
$ws.Range("F706").Value = 40676
$ws.Range("G706").Value = 4953
$ws.Range("F707").Value = 38819
$ws.Range("G707").Value = 4618
$ws.Range("F708").Value = 35602
$ws.Range("G708").Value = 4159
$ws.Range("F709").Value = 32386
$ws.Range("G709").Value = 3980
$ws.Range("F710").Value = 14702
$ws.Range("G710").Value = 2634
$ws.Range("F711").Value = 22568
$ws.Range("G711").Value = 3812
$ws.Range("F712").Value = 51321
$ws.Range("G712").Value = 6310
$ws.Range("F713").Value = 37135
$ws.Range("G713").Value = 4744
$ws.Range("F714").Value = 32363
$ws.Range("G714").Value = 3977
$ws.Range("F715").Value = 31698
$ws.Range("G715").Value = 3554
$ws.Range("F716").Value = 29646
$ws.Range("G716").Value = 3673
$ws.Range("F717").Value = 12464
$ws.Range("G717").Value = 2114
$ws.Range("F718").Value = 16916
$ws.Range("G718").Value = 2812
$ws.Range("F719").Value = 43393
$ws.Range("G719").Value = 5162
$ws.Range("F720").Value = 30980
$ws.Range("G720").Value = 3493
$ws.Range("F721").Value = 27790
$ws.Range("G721").Value = 3116
$ws.Range("F722").Value = 27790
$ws.Range("G722").Value = 2857
$ws.Range("F723").Value = 22143
$ws.Range("G723").Value = 2724
$ws.Range("F724").Value = 9310
$ws.Range("G724").Value = 1496
$ws.Range("F725").Value = 12550
$ws.Range("G725").Value = 2045
$ws.Range("F726").Value = 34752
$ws.Range("G726").Value = 4047
$ws.Range("F727").Value = 24708
$ws.Range("G727").Value = 2743
$ws.Range("F728").Value = 24359
$ws.Range("G728").Value = 2553

# --- Append new daily rows 729-732 (2022-03-03 .. 2022-03-06) ---
$ws.Range("A729").Value = 44623
$ws.Range("B729").Value = 1493383
$ws.Range("C729").Value = 19118
$ws.Range("D729").Value = 11029
$ws.Range("E729").Value = 18631
$ws.Range("F729").Value = 22181
$ws.Range("G729").Value = 2389
$ws.Range("A729").NumberFormat = "yyyy-mm-dd"

$ws.Range("A730").Value = 44624
$ws.Range("B730").Value = 1503308
$ws.Range("C730").Value = 18787
$ws.Range("D730").Value = 9925
$ws.Range("E730").Value = 18663
$ws.Range("F730").Value = 16061
$ws.Range("G730").Value = 1962
$ws.Range("A730").NumberFormat = "yyyy-mm-dd"

$ws.Range("A731").Value = 44625
$ws.Range("B731").Value = 1509494
$ws.Range("C731").Value = 11580
$ws.Range("D731").Value = 6186
$ws.Range("E731").Value = 18687
$ws.Range("F731").Value = 6472
$ws.Range("G731").Value = 1067
$ws.Range("A731").NumberFormat = "yyyy-mm-dd"

$ws.Range("A732").Value = 44626
$ws.Range("B732").Value = 1512913
$ws.Range("C732").Value = 6564
$ws.Range("D732").Value = 3419
$ws.Range("E732").Value = 18704
$ws.Range("F732").Value = 7086
$ws.Range("G732").Value = 1246
$ws.Range("A732").NumberFormat = "yyyy-mm-dd"

